# Edit script applying the 31/12/2025 16:48 update (LP1912+215+6203/6173)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 (19 new rows: 1106-1124) ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 16:48:51"
$ws1.Cells.Item(3, 1).Value = "Total filas: 1123"

$ws1.Cells.Item(1106, 2).Value = "16:48:41"
$ws1.Cells.Item(1106, 3).Value = "17:03"
$ws1.Cells.Item(1106, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(1106, 5).Value = 15
$ws1.Cells.Item(1106, 6).Value = "LP1912"
$ws1.Cells.Item(1106, 7).Value = "31/12/2025"
$ws1.Cells.Item(1107, 2).Value = "16:48:41"
$ws1.Cells.Item(1107, 3).Value = "17:04"
$ws1.Cells.Item(1107, 4).Value = "14_ABASTO"
$ws1.Cells.Item(1107, 5).Value = 16
$ws1.Cells.Item(1107, 6).Value = "LP1912"
$ws1.Cells.Item(1107, 7).Value = "31/12/2025"
$ws1.Cells.Item(1108, 2).Value = "16:48:41"
$ws1.Cells.Item(1108, 3).Value = "17:07"
$ws1.Cells.Item(1108, 4).Value = "15_ABASTO"
$ws1.Cells.Item(1108, 5).Value = 19
$ws1.Cells.Item(1108, 6).Value = "LP1912"
$ws1.Cells.Item(1108, 7).Value = "31/12/2025"
$ws1.Cells.Item(1109, 2).Value = "16:48:41"
$ws1.Cells.Item(1109, 3).Value = "17:14"
$ws1.Cells.Item(1109, 4).Value = "10_OLMOS"
$ws1.Cells.Item(1109, 5).Value = 26
$ws1.Cells.Item(1109, 6).Value = "LP1912"
$ws1.Cells.Item(1109, 7).Value = "31/12/2025"
$ws1.Cells.Item(1110, 2).Value = "16:48:41"
$ws1.Cells.Item(1110, 3).Value = "17:23"
$ws1.Cells.Item(1110, 4).Value = "16_SANTA ANA"
$ws1.Cells.Item(1110, 5).Value = 35
$ws1.Cells.Item(1110, 6).Value = "LP1912"
$ws1.Cells.Item(1110, 7).Value = "31/12/2025"
$ws1.Cells.Item(1111, 2).Value = "16:48:41"
$ws1.Cells.Item(1111, 3).Value = "17:24"
$ws1.Cells.Item(1111, 4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(1111, 5).Value = 36
$ws1.Cells.Item(1111, 6).Value = "LP1912"
$ws1.Cells.Item(1111, 7).Value = "31/12/2025"
$ws1.Cells.Item(1112, 2).Value = "16:48:41"
$ws1.Cells.Item(1112, 3).Value = "17:27"
$ws1.Cells.Item(1112, 4).Value = "15_ABASTO"
$ws1.Cells.Item(1112, 5).Value = 39
$ws1.Cells.Item(1112, 6).Value = "LP1912"
$ws1.Cells.Item(1112, 7).Value = "31/12/2025"
$ws1.Cells.Item(1113, 2).Value = "16:48:41"
$ws1.Cells.Item(1113, 3).Value = "17:33"
$ws1.Cells.Item(1113, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(1113, 5).Value = 45
$ws1.Cells.Item(1113, 6).Value = "LP1912"
$ws1.Cells.Item(1113, 7).Value = "31/12/2025"
$ws1.Cells.Item(1114, 2).Value = "16:48:41"
$ws1.Cells.Item(1114, 3).Value = "17:34"
$ws1.Cells.Item(1114, 4).Value = "10_OLMOS"
$ws1.Cells.Item(1114, 5).Value = 46
$ws1.Cells.Item(1114, 6).Value = "LP1912"
$ws1.Cells.Item(1114, 7).Value = "31/12/2025"
$ws1.Cells.Item(1115, 2).Value = "16:48:41"
$ws1.Cells.Item(1115, 3).Value = "17:35"
$ws1.Cells.Item(1115, 4).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(1115, 5).Value = 47
$ws1.Cells.Item(1115, 6).Value = "LP1912"
$ws1.Cells.Item(1115, 7).Value = "31/12/2025"
$ws1.Cells.Item(1116, 2).Value = "16:48:41"
$ws1.Cells.Item(1116, 3).Value = "17:38"
$ws1.Cells.Item(1116, 4).Value = "17X38_ROMERO"
$ws1.Cells.Item(1116, 5).Value = 50
$ws1.Cells.Item(1116, 6).Value = "LP1912"
$ws1.Cells.Item(1116, 7).Value = "31/12/2025"
$ws1.Cells.Item(1117, 2).Value = "16:48:41"
$ws1.Cells.Item(1117, 3).Value = "17:47"
$ws1.Cells.Item(1117, 4).Value = "16_SANTA ANA"
$ws1.Cells.Item(1117, 5).Value = 59
$ws1.Cells.Item(1117, 6).Value = "LP1912"
$ws1.Cells.Item(1117, 7).Value = "31/12/2025"
$ws1.Cells.Item(1118, 2).Value = "16:48:41"
$ws1.Cells.Item(1118, 3).Value = "17:50"
$ws1.Cells.Item(1118, 4).Value = "215_EL PELIGRO"
$ws1.Cells.Item(1118, 5).Value = 62
$ws1.Cells.Item(1118, 6).Value = "LP1912"
$ws1.Cells.Item(1118, 7).Value = "31/12/2025"
$ws1.Cells.Item(1119, 2).Value = "16:48:41"
$ws1.Cells.Item(1119, 3).Value = "17:54"
$ws1.Cells.Item(1119, 4).Value = "10_OLMOS"
$ws1.Cells.Item(1119, 5).Value = 66
$ws1.Cells.Item(1119, 6).Value = "LP1912"
$ws1.Cells.Item(1119, 7).Value = "31/12/2025"
$ws1.Cells.Item(1120, 2).Value = "16:48:41"
$ws1.Cells.Item(1120, 3).Value = "17:59"
$ws1.Cells.Item(1120, 4).Value = "16_SANTA ANA"
$ws1.Cells.Item(1120, 5).Value = 71
$ws1.Cells.Item(1120, 6).Value = "LP1912"
$ws1.Cells.Item(1120, 7).Value = "31/12/2025"
$ws1.Cells.Item(1121, 2).Value = "16:48:41"
$ws1.Cells.Item(1121, 3).Value = "18:02"
$ws1.Cells.Item(1121, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(1121, 5).Value = 74
$ws1.Cells.Item(1121, 6).Value = "LP1912"
$ws1.Cells.Item(1121, 7).Value = "31/12/2025"
$ws1.Cells.Item(1122, 2).Value = "16:48:41"
$ws1.Cells.Item(1122, 3).Value = "18:04"
$ws1.Cells.Item(1122, 4).Value = "14_ABASTO"
$ws1.Cells.Item(1122, 5).Value = 76
$ws1.Cells.Item(1122, 6).Value = "LP1912"
$ws1.Cells.Item(1122, 7).Value = "31/12/2025"
$ws1.Cells.Item(1123, 2).Value = "16:48:41"
$ws1.Cells.Item(1123, 3).Value = "18:24"
$ws1.Cells.Item(1123, 4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(1123, 5).Value = 96
$ws1.Cells.Item(1123, 6).Value = "LP1912"
$ws1.Cells.Item(1123, 7).Value = "31/12/2025"
$ws1.Cells.Item(1124, 2).Value = "16:48:41"
$ws1.Cells.Item(1124, 3).Value = "18:27"
$ws1.Cells.Item(1124, 4).Value = "15_ABASTO"
$ws1.Cells.Item(1124, 5).Value = 99
$ws1.Cells.Item(1124, 6).Value = "LP1912"
$ws1.Cells.Item(1124, 7).Value = "31/12/2025"

# --- Sheet 2: LP1912-215 (1 new row: 76) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 16:48:51"
$ws2.Cells.Item(3, 1).Value = "Total filas: 75"
$ws2.Cells.Item(76, 2).Value = "31/12/2025"
$ws2.Cells.Item(76, 3).Value = "16:48:41"
$ws2.Cells.Item(76, 4).Value = "17:50"
$ws2.Cells.Item(76, 5).Value = "215_EL PELIGRO"
$ws2.Cells.Item(76, 6).Value = 62
$ws2.Cells.Item(76, 7).Value = "LP1912"

# --- Sheet 3: 6203-6173 (1 new row: 137) ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 16:48:51"
$ws3.Cells.Item(3, 1).Value = "Total filas: 136"
$ws3.Cells.Item(137, 2).Value = "31/12/2025"
$ws3.Cells.Item(137, 3).Value = "16:48:46"
$ws3.Cells.Item(137, 4).Value = "16:57"
$ws3.Cells.Item(137, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(137, 6).Value = 9
$ws3.Cells.Item(137, 7).Value = "L6203"

